$d = $word.ActiveDocument

function Get-ParaIndexByText($text) {
    $n = $d.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]7,[char]12,[char]13)
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

# --- Step 1: rename "Discrete Math" -> "Digital Logic" ---
$idx = Get-ParaIndexByText("Discrete Math")
$p = $d.Paragraphs.Item($idx)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$r.Text = "Digital Logic"

# --- Step 2: rename "Object Oriented Programming" -> "Discrete Math" ---
$idx = Get-ParaIndexByText("Object Oriented Programming")
$p = $d.Paragraphs.Item($idx)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$r.Text = "Discrete Math"

# --- Step 3: insert the new checklist items after the (renamed) "Discrete Math" paragraph ---
$newItems = @("Distributed Computing", "Embedded Systems", "Graphics", "Human-Computer Interaction", "Networks", "Object Oriented Programming")

$idx = Get-ParaIndexByText("Discrete Math")
$anchor = $d.Paragraphs.Item($idx)
$ar = $anchor.Range
$ar.Collapse(0)

foreach ($item in $newItems) {
    $ar.InsertAfter($item + [char]13)
    $newIdx = Get-ParaIndexByText($item)
    $ar = $d.Paragraphs.Item($newIdx).Range
    $ar.Collapse(0)
}

# --- Step 4: split "Operating Systems" off from its trailing page break, ---
# --- then add "Software Engineering" (with the page break) right after  ---
$idx = Get-ParaIndexByText("Operating Systems")
$p = $d.Paragraphs.Item($idx)
$rng = $p.Range
$splitPos = $rng.End - 2
$rSplit = $d.Range($splitPos, $splitPos)
$rSplit.InsertAfter([char]13)

$idx = Get-ParaIndexByText("Operating Systems")
$p = $d.Paragraphs.Item($idx)
$r2 = $p.Range
$r2.Collapse(0)
$r2.InsertAfter("Software Engineering" + [char]13)

$idx = Get-ParaIndexByText("Software Engineering")
$p = $d.Paragraphs.Item($idx)
$markPos = $p.Range.End - 1
$markRange = $d.Range($markPos, $markPos + 1)
$markRange.Delete()

# --- Step 5: move the "_GoBack" bookmark so it sits right after "Graphics" ---
$idx = Get-ParaIndexByText("Graphics")
$p = $d.Paragraphs.Item($idx)
$endPos = $p.Range.End - 1

$rEnd = $d.Range($endPos, $endPos)
$rEnd.InsertAfter("X")

$rb = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $rb)

$rDel = $d.Range($endPos, $endPos + 1)
$rDel.Delete()

Write-Output "done"
